$d = $word.ActiveDocument

# The document contains seven occurrences of a split "<id>...</id>" marker,
# each built from three separate runs:
#   run1 (Courier New, color 7f6000): "<id>"
#   run2 (color 000000):              "p046v_N"
#   run3 (Courier New, color 7f6000): "</id>"
# Collapse each triple into a single run (keeping run1's formatting) whose
# text is the concatenation "<id>p046v_N</id>".
for ($i = 1; $i -le 7; $i++) {
    $needle = "<id>p046v_$i</id>"
    $d.Content.Find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $needle, 2) | Out-Null
}
